{"js": "// Word JavaScript API (Office.js) script\n// Body of: async (context) => { ... }\n//\n// Changes implemented (per the authoritative XML diff):\n// 1. In the \"Skills & Abilities\" table, after the paragraph that ends with\n//    \"I have some experience with the MacOS system as well.\" add three new\n//    paragraphs:\n//      - \"I help to care for my cousin with Down\\u2019s Syndrome at home.\"\n//      - \"I am CPR/AED/First Aid Certified.\"\n//      - \"I have been trained in the PA Modified Medication Administration\n//         Training Course. This includes the administration of insulin.\"\n// 2. In the Education table, change \"...Accounting Minor, Pittsburgh...\" to\n//    \"...Accounting Major, Pittsburgh...\" (i.e. Minor -> Major).\n\nconst body = context.document.body;\n\n// --- Change 1: locate the MacOS paragraph and append the three new ones ---\nconst macResults = body.search(\n  \"I have some experience with the MacOS system as well.\",\n  { matchCase: true, matchWholeWord: false }\n);\nmacResults.load(\"items\");\nawait context.sync();\n\nif (macResults.items.length === 0) {\n  throw new Error(\"Could not find the MacOS skills paragraph to anchor the insert.\");\n}\n\nconst anchor = macResults.items[0];\n\n// Insert the three new paragraphs, each right after the anchor range, so the\n// final reading order is: cousin / CPR / medication-training.\nconst medTrainingPara = anchor.insertParagraph(\n  \"I have been trained in the PA Modified Medication Administration Training Course.\",\n  Word.InsertLocation.after\n);\nmedTrainingPara.insertText(\" This includes the administration of insulin.\", Word.InsertLocation.end);\n\nanchor.insertParagraph(\"I am CPR/AED/First Aid Certified.\", Word.InsertLocation.after);\n\nanchor.insertParagraph(\n  \"I help to care for my cousin with Down\\u2019s Syndrome at home.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n\n// --- Change 2: \"Minor\" -> \"Major\" in the Education entry ---\nconst minorResults = body.search(\"Accounting Minor, Pittsburgh, Point Park University\", { matchCase: true });\nminorResults.load(\"items\");\nawait context.sync();\n\nif (minorResults.items.length === 0) {\n  throw new Error(\"Could not find the 'Accounting Minor, Pittsburgh, Point Park University' text.\");\n}\n\nminorResults.items[0].insertText(\n  \"Accounting Major, Pittsburgh, Point Park University\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script\n# $word.ActiveDocument / $d is the open document.\n#\n# Changes implemented (per the authoritative XML diff):\n# 1. In the \"Skills & Abilities\" table, after the paragraph that ends with\n#    \"I have some experience with the MacOS system as well.\" add three new\n#    paragraphs:\n#      - \"I help to care for my cousin with Down\u2019s Syndrome at home.\"\n#      - \"I am CPR/AED/First Aid Certified.\"\n#      - \"I have been trained in the PA Modified Medication Administration\n#         Training Course. This includes the administration of insulin.\"\n# 2. In the Education table, change \"...Accounting Minor, Pittsburgh...\" to\n#    \"...Accounting Major, Pittsburgh...\" (i.e. Minor -> Major).\n\n$d = $word.ActiveDocument\n\n# --- Change 1: find the MacOS paragraph and append the three new ones ---\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*I have some experience with the MacOS system as well.*\") {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -eq $null) {\n    throw \"Could not find the MacOS skills paragraph to anchor the insert.\"\n}\n\n# Paragraph 1: cousin / Down's Syndrome\n$anchor.Range.InsertParagraphAfter()\n$p1 = $anchor.Next()\n$p1.Range.Text = \"I help to care for my cousin with Down\u2019s Syndrome at home.\"\n\n# Paragraph 2: CPR/AED/First Aid\n$p1.Range.InsertParagraphAfter()\n$p2 = $p1.Next()\n$p2.Range.Text = \"I am CPR/AED/First Aid Certified.\"\n\n# Paragraph 3: PA Modified Medication Administration Training Course\n$p2.Range.InsertParagraphAfter()\n$p3 = $p2.Next()\n$p3.Range.Text = \"I have been trained in the PA Modified Medication Administration Training Course. This includes the administration of insulin.\"\n\n# --- Change 2: \"Minor\" -> \"Major\" in the Education entry ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Accounting Minor, Pittsburgh, Point Park University\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"Accounting Major, Pittsburgh, Point Park University\"\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$find.Execute($null, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $null, $wdReplaceAll)\n\nWrite-Output \"done\"\n"}
